$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("N2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("T2").Value = 0

# Row 4 values
$ws.Range("N4").Value = 0.5607188170993559
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0.1214376341987118
$ws.Range("Q4").Value = 0.7831777282599498
$ws.Range("S4").Value = 0.5663554565198996
$ws.Range("T4").Value = 0.7087799791449426
$ws.Range("U4").Value = -1
$ws.Range("V4").Value = 0.4175599582898852

# Column B values rows 15-21
$ws.Range("B15").Value = 9.921850020943705
$ws.Range("B16").Value = 5.021344692993106
$ws.Range("B17").Value = -8.327466958658352
$ws.Range("B18").Value = 6.707646015967205
$ws.Range("B19").Value = 9.875561065351537
$ws.Range("B20").Value = 13.11182579371309
$ws.Range("B21").Value = 25.60214361037885
